# This script applies a 3-way rotation of the data rows 56/58/59 and a
# 2-way swap of specific fields in rows 60/61 in the "Artfynd" sheet,
# as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 56  (becomes what used to be in row 58: "Blåsippa" record)
# ---------------------------------------------------------------
$ws.Range("A56").Value = 111898889
$ws.Range("B56").Value = 98535
$ws.Range("D56").Value = "LC"
$ws.Range("E56").Value = 222498
$ws.Range("F56").Value = "Blåsippa"
$ws.Range("G56").Value = "Hepatica nobilis"
$ws.Range("H56").Value = "Schreb."
$ws.Range("K56").Value = "fullt utvecklade blad"
$ws.Range("Q56").Value = 650135.0421630922
$ws.Range("R56").Value = 6654002.501842719
$ws.Range("AC56").ClearContents()
$ws.Range("AH56").Value = "Ängsbarrskog"
$ws.Range("AI56").Value = "Ungskog"

# ---------------------------------------------------------------
# Row 58  (becomes what used to be in row 59: "Svavelriska" record)
# ---------------------------------------------------------------
$ws.Range("A58").Value = 111898191
$ws.Range("B58").Value = 90332
$ws.Range("E58").Value = 4769
$ws.Range("F58").Value = "Svavelriska"
$ws.Range("G58").Value = "Lactarius scrobiculatus"
$ws.Range("H58").Value = "(Scop.:Fr.) Fr."
$ws.Range("I58").NumberFormat = "@"
$ws.Range("I58").Value = "2"
$ws.Range("J58").Value = "fruktkroppar"
$ws.Range("K58").ClearContents()

# ---------------------------------------------------------------
# Row 59  (becomes what used to be in row 56: "Skogsalm" record)
# ---------------------------------------------------------------
$ws.Range("A59").Value = 111898660
$ws.Range("B59").Value = 100532
$ws.Range("D59").Value = "CR"
$ws.Range("E59").Value = 223246
$ws.Range("F59").Value = "Skogsalm"
$ws.Range("G59").Value = "Ulmus glabra"
$ws.Range("H59").Value = "Huds."
$ws.Range("I59").ClearContents()
$ws.Range("J59").ClearContents()
$ws.Range("Q59").Value = 650054.1336129439
$ws.Range("R59").Value = 6654018.240072312
$ws.Range("AC59").Value = "Stammens omkrets i brösthöjd: 64 cm"
$ws.Range("AH59").Value = "Ängsblandskog"
$ws.Range("AI59").ClearContents()

# ---------------------------------------------------------------
# Rows 60/61 (Knärot records): swap Id, Antal, Ost, Nord
# ---------------------------------------------------------------
$ws.Range("A60").Value = 111911660
$ws.Range("I60").NumberFormat = "@"
$ws.Range("I60").Value = "19"
$ws.Range("Q60").Value = 650026.652882754
$ws.Range("R60").Value = 6654299.07778531

$ws.Range("A61").Value = 111911698
$ws.Range("I61").NumberFormat = "@"
$ws.Range("I61").Value = "16"
$ws.Range("Q61").Value = 650032.9755174413
$ws.Range("R61").Value = 6654279.303373625
